$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data grid (rows 2..10), cross product of ECs/FAPs/MuSCs sending x target clusters,
# with Myoc/Fzd4 ligand/receptor, and updated TPM-derived statistics.
$data = @(
    @("ECs",   "Myoc", "Fzd4", "ECs",   3, 1, 0.720897,            2.162691,            0.0284720950782092,  0.02847209507820921, 3, 1, 34.97976933333334, 104.939308,          0.4352965780925344, 0.4352965780925344, 25.216810773092,    226.951296957828,   0.01239380555866976,  0.01239380555866976),
    @("ECs",   "Myoc", "Fzd4", "FAPs",  3, 1, 0.720897,            2.162691,            0.0284720950782092,  0.02847209507820921, 3, 1, 20.343383,          61.03014900000001,   0.2531579017099818, 0.2531579017099818, 14.665483774551,    131.989353970959,   0.007207935847286543, 0.007207935847286543),
    @("ECs",   "Myoc", "Fzd4", "MuSCs", 3, 1, 0.720897,            2.162691,            0.0284720950782092,  0.02847209507820921, 3, 1, 25.035323,          75.105969,           0.3115455201974837, 0.3115455201974837, 18.047889244731,    162.431003202579,   0.008870353672252902, 0.008870353672252904),
    @("FAPs",  "Myoc", "Fzd4", "ECs",   3, 1, 23.971258,           71.91377399999999,  0.9467537483444692,  0.9467537483444693,  3, 1, 34.97976933333334, 104.939308,          0.4352965780925344, 0.4352965780925344, 838.5090754698213,  7546.581679228391,  0.4121186669506279,   0.4121186669506279),
    @("FAPs",  "Myoc", "Fzd4", "FAPs",  3, 1, 23.971258,           71.91377399999999,  0.9467537483444692,  0.9467537483444693,  3, 1, 20.343383,          61.03014900000001,   0.2531579017099818, 0.2531579017099818, 487.656482485814,   4388.908342372326,  0.239678192366946,    0.239678192366946),
    @("FAPs",  "Myoc", "Fzd4", "MuSCs", 3, 1, 23.971258,           71.91377399999999,  0.9467537483444692,  0.9467537483444693,  3, 1, 25.035323,          75.105969,           0.3115455201974837, 0.3115455201974837, 600.1281867463339,  5401.153680717005,  0.2949568890268953,   0.2949568890268953),
    @("MuSCs", "Myoc", "Fzd4", "ECs",   3, 1, 0.6272673333333333, 1.881802,            0.0247741565773216,  0.0247741565773216,  3, 1, 34.97976933333334, 104.939308,          0.4352965780925344, 0.4352965780925344, 21.94166663033511,  197.474999673016,   0.01078410558323675,  0.01078410558323675),
    @("MuSCs", "Myoc", "Fzd4", "FAPs",  3, 1, 0.6272673333333333, 1.881802,            0.0247741565773216,  0.0247741565773216,  3, 1, 20.343383,          61.03014900000001,   0.2531579017099818, 0.2531579017099818, 12.76073960538867,  114.846656448498,   0.006271773495749281, 0.006271773495749281),
    @("MuSCs", "Myoc", "Fzd4", "MuSCs", 3, 1, 0.6272673333333333, 1.881802,            0.0247741565773216,  0.0247741565773216,  3, 1, 25.035323,          75.105969,           0.3115455201974837, 0.3115455201974837, 15.70384029734867,  141.334562676138,   0.00771827749833557,  0.007718277498335572)
)

# Remove the old rows 2..13 (12 data rows) and rewrite with the new 9 data rows (rows 2..10)
$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count
if ($lastRow -gt 1) {
    $ws.Range("A2:T$lastRow").ClearContents()
}

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $values = $data[$i]
    for ($c = 0; $c -lt $values.Count; $c++) {
        $col = $c + 1
        $ws.Cells.Item($row, $col).Value = $values[$c]
    }
}

$excel.ActiveSheet.UsedRange | Out-Null
